$d = $word.ActiveDocument

# The first row/first column cell of the table currently contains a single
# empty paragraph. Insert a brand-new paragraph with the text
# "La mona pelona" right before that existing (empty) paragraph, leaving
# the empty paragraph in place after it.
$cell = $d.Tables(1).Cell(1, 1)
$cell.Range.InsertParagraphBefore() | Out-Null

# The insertion is a structural edit, so re-fetch the cell/paragraph handles
# before touching them again. The new paragraph already inherits the
# es-ES language run formatting from the paragraph mark it was split from.
$newPara = $d.Tables(1).Cell(1, 1).Range.Paragraphs(1).Range
$newPara.InsertBefore("La mona pelona")
